$wb = $excel.ActiveWorkbook

# Insert the new "Week 11" sheet. With no args, Add() inserts before the
# active sheet ("Week 10", the current first/active tab) and makes the new
# sheet active/selected - matching how the other week tabs get pushed back.
$ws = $wb.Worksheets.Add()
$ws.Name = "Week 11"

$data = @(
    @("game", "temp", "wind"),
    @("NEvsNYJ", 39, 8),
    @("MIAvsWAS", 58, 10),
    @("ATLvsCAR", 63, 8),
    @("BUFvsTB", 43, 9),
    @("HOUvsTEN", 66, 11),
    @("CHIvsMIN", 47, 9),
    @("GBvsNYG", 45, 9),
    @("CINvsPIT", 60, 12),
    @("JAXvsLAC", 76, 8),
    @("LAvsSEA", 63, 7),
    @("ARIvsSF", 67, 3),
    @("BALvsCLE", 56, 11),
    @("DENvsKC", 54, 4),
    @("DETvsPHI", 55, 10),
    @("DALvsLV", 56, 3)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 1
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}

$ws.Range("C16").Select() | Out-Null
